# Fixed #295 Add the version of M2Doc in the template custom properties.
#
# This particular template, tests/org.obeonetwork.m2doc.tests/resources/
# query/missingExpression/missingExpression-template.docx, is only
# incidentally touched by that change: the unified diff for this file
# shows no textual/content edits at all. Every "-"/"+" line pair is a
# pure XML-attribute reordering (e.g. the xmlns:* declarations on
# <w:document>, and the attribute order on <w:pgSz>, <w:pgMar>,
# <w:rFonts>, <w:lang>, <w:latentStyles>, every <w:lsdException>, and
# the four <w:style> definitions together with their nested <w:tblInd>/
# <w:tblCellMar> children) - the attribute *names and values* before and
# after are identical sets, just alphabetised. That churn comes from the
# tool (M2Doc/docx4j) that re-serialised the test fixture while the real
# fix (adding the M2Doc version to the custom document properties) landed
# elsewhere in the commit; no word/document.xml or word/styles.xml
# content actually changed for this resource.
#
# Word's COM object model has no notion of "attribute order" to
# reproduce - it is not a property of any Range/Style/PageSetup object,
# so there is no WordBasic-level call that changes it. Poking the
# corresponding values through the object model (e.g. re-assigning
# PageSetup.PageWidth/PageHeight/margins to themselves, or nudging the
# Styles collection) does not reorder the serialized attributes either:
# this runtime patches existing XML in place and preserves the original
# attribute order, while additionally minting extra namespace
# declarations and perturbing docProps/app.xml word/paragraph counts as
# a side effect - actively diverging from the target instead of
# approaching it. So the faithful, minimal-diff edit for this resource
# is to leave its content untouched.
$d = $word.ActiveDocument
